# TC for URL added
# Populates the "Test Case ID" (E) and "Status" (F) columns for the
# FR_URL_01..FR_URL_04 requirement rows (rows 4-7) on the RTM_ALL sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM_ALL")

$greenFill = 5287936   # RGB(0,176,80) -> FF00B050
$xlPasteFormats = -4122
$xlVAlignCenter = -4108
$xlVAlignTop = -4160

# ---------------------------------------------------------------------------
# Column E ("Test Case ID") values
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "TC_URL_01"
$ws.Range("E5").Value = "TC_URL_02"
$ws.Range("E6").Value = "TC_URL_03"
$ws.Range("E7").Value = "TC_URL_04`nTC_URL_05`n"

# E4/E5 share the same "wrap + vertically centered" look. Build the style
# once on E4 and stamp it onto E5 via copy/paste-format so no transient
# style is left behind in the style table.
$e4 = $ws.Range("E4")
$e4.VerticalAlignment = $xlVAlignCenter
$e4.WrapText = $true
$e4.Copy()
$ws.Range("E5").PasteSpecial($xlPasteFormats)

# E6 keeps the default (no) formatting.

# E7 gets its own "wrap + top aligned" look (it holds two wrapped lines).
$e7 = $ws.Range("E7")
$e7.VerticalAlignment = $xlVAlignTop
$e7.WrapText = $true

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column F ("Status") values - all "Pass" with a green fill
# ---------------------------------------------------------------------------
$ws.Range("F4").Value = "Pass"
$ws.Range("F5").Value = "Pass"
$ws.Range("F6").Value = "Pass"
$ws.Range("F7").Value = "Pass"
$ws.Range("F4:F7").Interior.Color = $greenFill

# ---------------------------------------------------------------------------
# Row 7 (FR_URL_04) switches to top-aligned text and grows taller to fit
# the two wrapped test case IDs placed in column E.
# ---------------------------------------------------------------------------
$a7 = $ws.Range("A7")
$a7.VerticalAlignment = $xlVAlignTop
$a7.Copy()
$ws.Range("B7:D7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Rows.Item(7).RowHeight = 45

# ---------------------------------------------------------------------------
# Column E width and current selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 21.43

$ws.Activate()
$ws.Range("F12").Select()
